$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B4: "consult statistician" -> "planning with experts, consult statistician"
$ws.Range("B4").Value = "planning with experts, consult statistician"

# Delete row 10 ("plan" / "planning with experts") entirely; rows below shift up.
$ws.Rows.Item(10).Delete()
